# Update cryptocurrency price (D) and volume change (E) values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.656.28'
$ws.Range("E2").Value = '  -2.97%  '

$ws.Range("D3").Value = '2.894.05'
$ws.Range("E3").Value = '  -4.19%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '584.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.22%  '

$ws.Range("E6").Value = '  -2.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.503'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -3.23%  '

$ws.Range("D9").Value = '2.894.24'
$ws.Range("E9").Value = '  -4.16%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +4.64%  '

$ws.Range("E11").Value = '  -4.83%  '

$ws.Range("E12").Value = '  -3.10%  '

$ws.Range("E13").Value = '  -4.31%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '33.77'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.14%  '

$ws.Range("E15").Value = '  +0.48%  '

$ws.Range("D16").Value = '3.372.15'
$ws.Range("E16").Value = '  -4.22%  '

$ws.Range("D17").Value = '60.562.01'
$ws.Range("E17").Value = '  -2.98%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.78'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.49%  '

$ws.Range("D19").Value = '2.894.99'
$ws.Range("E19").Value = '  -3.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '424.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.41%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.57'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.53%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.666'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.57%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -5.40%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '79.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.87%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.17%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.20'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.98%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.82'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.99%  '

$ws.Range("E28").Value = '  -0.03%  '

$ws.Range("E29").Value = '  +0.08%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.19'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.90%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.18'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.28%  '

$ws.Range("E32").Value = '  -3.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '26.29'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -4.74%  '

$ws.Range("E34").Value = '  -3.96%  '

$ws.Range("D35").Value = '0.0₃0829'
$ws.Range("E35").Value = '  -2.64%  '

$ws.Range("E36").Value = '  -2.74%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.91%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '49.24'
$ws.Range("D38").Style = "Normal"

$ws.Range("E39").Value = '  -2.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.92'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.35%  '

$ws.Range("E41").Value = '  +0.06%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.67'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.13%  '

$ws.Range("E43").Value = '  +1.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.94%  '

$ws.Range("E45").Value = '  -2.69%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '370.77'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -5.35%  '

$ws.Range("D47").Value = '2.646.63'
$ws.Range("E47").Value = '  -3.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.95'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.19%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '24.80'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.34%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.106'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.74%  '
